$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E535").Value = 0.04721407624633431
$ws.Range("I535").Value = 0.008379839814381316
$ws.Range("K535").Value = 0.8379839814381316
$ws.Range("L535").Value = 3.616597777396834
$ws.Range("E536").Value = 0.04877049180327869
$ws.Range("G536").Value = 0.04721407624633431
$ws.Range("I536").Value = 0.002844330519725851
$ws.Range("K536").Value = 0.2844330519725851
$ws.Range("L536").Value = 11.59383894276562
$ws.Range("E537").Value = 0.05651058370750481
$ws.Range("G537").Value = 0.04877049180327869
$ws.Range("I537").Value = 0.001953306137813373
$ws.Range("K537").Value = 0.1953306137813373
$ws.Range("L537").Value = 20.6202834197878
$ws.Range("E538").Value = 0.05929203539823009
$ws.Range("G538").Value = 0.05651058370750481
$ws.Range("I538").Value = 0.01072348407663616
$ws.Range("K538").Value = 1.072348407663616
$ws.Range("L538").Value = 2.37508410014622
$ws.Range("E539").Value = 0.05481481481481482
$ws.Range("G539").Value = 0.05929203539823009
$ws.Range("I539").Value = 0.003261720761925913
$ws.Range("K539").Value = 0.3261720761925913
$ws.Range("L539").Value = 9.532281683442914
$ws.Range("E540").Value = 0.05855491329479768
$ws.Range("G540").Value = 0.05481481481481482
$ws.Range("I540").Value = 0.001899283610588542
$ws.Range("K540").Value = 0.1899283610588542
$ws.Range("L540").Value = 23.53346872896167
$ws.Range("E541").Value = 0.04393063583815029
$ws.Range("G541").Value = 0.05855491329479768
$ws.Range("I541").Value = 0.007845768716951331
$ws.Range("K541").Value = 0.7845768716951331
$ws.Range("L541").Value = 3.658552831029914
$ws.Range("E542").Value = 0.04724919093851133
$ws.Range("G542").Value = 0.04393063583815029
$ws.Range("I542").Value = 0.002822284678592656
$ws.Range("K542").Value = 0.2822284678592656
$ws.Range("L542").Value = 10.88973531777122
$ws.Range("E543").Value = 0.05115712545676004
$ws.Range("G543").Value = 0.04724919093851133
$ws.Range("I543").Value = 0.001805175545496719
$ws.Range("K543").Value = 0.1805175545496719
$ws.Range("L543").Value = 21.51915997808035
$ws.Range("E544").Value = 0.05527950310559006
$ws.Range("G544").Value = 0.05115712545676004
$ws.Range("I544").Value = 0.009697904295128781
$ws.Range("K544").Value = 0.9697904295128781
$ws.Range("L544").Value = 3.434376476767386
$ws.Range("E545").Value = 0.05106382978723405
$ws.Range("G545").Value = 0.05527950310559006
$ws.Range("I545").Value = 0.003015557014504368
$ws.Range("K545").Value = 0.3015557014504368
$ws.Range("L545").Value = 10.54630634733379
$ws.Range("E546").Value = 0.05381381381381382
$ws.Range("G546").Value = 0.05106382978723405
$ws.Range("I546").Value = 0.001846291847303443
$ws.Range("K546").Value = 0.1846291847303443
$ws.Range("L546").Value = 22.09524879896936
$ws.Range("E547").Value = 0.04972067039106145
$ws.Range("G547").Value = 0.05381381381381382
$ws.Range("I547").Value = 0.008722696041987339
$ws.Range("K547").Value = 0.8722696041987339
$ws.Range("L547").Value = 3.818344033176162
$ws.Range("E548").Value = 0.06968253968253968
$ws.Range("G548").Value = 0.04972067039106145
$ws.Range("I548").Value = 0.004159620056971032
$ws.Range("K548").Value = 0.4159620056971032
$ws.Range("L548").Value = 7.402994890180952
$ws.Range("E549").Value = 0.06797004991680532
$ws.Range("G549").Value = 0.06968253968253968
$ws.Range("I549").Value = 0.002426902683360134
$ws.Range("K549").Value = 0.2426902683360134
$ws.Range("L549").Value = 15.66996034835925
$ws.Range("E550").Value = 0.05934718100890208
$ws.Range("G550").Value = 0.06797004991680532
$ws.Range("I550").Value = 0.01025852907355735
$ws.Range("K550").Value = 1.025852907355735
$ws.Range("L550").Value = 3.620137845334396
$ws.Range("E551").Value = 0.06022598870056498
$ws.Range("G551").Value = 0.05934718100890208
$ws.Range("I551").Value = 0.003392519742368263
$ws.Range("K551").Value = 0.3392519742368263
$ws.Range("L551").Value = 10.68337602154432
$ws.Range("E552").Value = 0.07699186991869918
$ws.Range("G552").Value = 0.06022598870056498
$ws.Range("I552").Value = 0.002576583998595741
$ws.Range("K552").Value = 0.2576583998595741
$ws.Range("L552").Value = 16.50337981705743
$ws.Range("E553").Value = 0.04651639344262296
$ws.Range("G553").Value = 0.07699186991869918
$ws.Range("I553").Value = 0.007898786598893148
$ws.Range("K553").Value = 0.7898786598893148
$ws.Range("L553").Value = 5.287262901623417
$ws.Range("E554").Value = 0.04721065778517902
$ws.Range("G554").Value = 0.04651639344262296
$ws.Range("I554").Value = 0.002606426956977766
$ws.Range("K554").Value = 0.2606426956977766
$ws.Range("L554").Value = 14.63289851640114
$ws.Range("E555").Value = 0.08487752928647498
$ws.Range("G555").Value = 0.04721065778517902
$ws.Range("I555").Value = 0.003062202623166073
$ws.Range("K555").Value = 0.3062202623166073
$ws.Range("L555").Value = 12.18457856988934
$ws.Range("E556").Value = 0.06577380952380953
$ws.Range("G556").Value = 0.08487752928647498
$ws.Range("I556").Value = 0.01121270666645249
$ws.Range("K556").Value = 1.121270666645249
$ws.Range("L556").Value = 3.633567235440834
$ws.Range("E557").Value = 0.05972944849115504
$ws.Range("G557").Value = 0.06577380952380953
$ws.Range("I557").Value = 0.003284124338305486
$ws.Range("K557").Value = 0.3284124338305486
$ws.Range("L557").Value = 11.73075294711453
$ws.Range("E558").Value = 0.06405188387893761
$ws.Range("G558").Value = 0.05972944849115504
$ws.Range("I558").Value = 0.002054626263777158
$ws.Range("K558").Value = 0.2054626263777158
$ws.Range("L558").Value = 22.13097215609249
